$d = $word.ActiveDocument

# 1. Update the two cached DATE field results: 15.01.2020 -> 17.01.2020
$range = $d.Content
$range.Find.Execute("15.01.2020", $false, $false, $false, $false, $false, $true, 1, $false, "17.01.2020", 2)
